$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.24418859991079245
$ws.Range("B1").Value = 0.24399259722358835
$ws.Range("A2").Value = -0.22188726358183253
$ws.Range("B2").Value = 0.22108492491182918
$ws.Range("A3").Value = -0.098779060104366323
$ws.Range("B3").Value = 0.098507057833987943
$ws.Range("A4").Value = -0.09050705786598634
$ws.Range("B4").Value = 0.090063697668400167
$ws.Range("A5").Value = -0.087063697687487895
$ws.Range("B5").Value = 0.085559366065469433
$ws.Range("A6").Value = -0.0031664588063495813
$ws.Range("B6").Value = 0.0030796630020866189
$ws.Range("A7").Value = 0.0069203369539381221
$ws.Range("B7").Value = -0.0069315191452976599
$ws.Range("A8").Value = 0.016931519101543113
$ws.Range("B8").Value = -0.016947188482680176
$ws.Range("A9").Value = 0.018947188463324771
$ws.Range("B9").Value = -0.018964448666509703
$ws.Range("A10").Value = 0.048271154475017397
$ws.Range("B10").Value = -0.048284295617625261
$ws.Range("A11").Value = 0.051284295596730978
$ws.Range("B11").Value = -0.051315246821059368
$ws.Range("A12").Value = -0.0079217571716823088
$ws.Range("B12").Value = 0.0078865632319100598
$ws.Range("A13").Value = -0.004386563253961917
$ws.Range("B13").Value = 0.0043766379193259652
$ws.Range("A14").Value = 0.003623362045256151
$ws.Range("B14").Value = -0.003623683997929561
$ws.Range("A15").Value = -0.0080510246295206755
$ws.Range("B15").Value = 0.0080331244245721578
$ws.Range("A16").Value = -0.0060331244425957387
$ws.Range("B16").Value = 0.006003449434534236
$ws.Range("A17").Value = -0.0040034494529974651
$ws.Range("B17").Value = 0.0039999999755169213
$ws.Range("A18").Value = -0.01610533366048017
$ws.Range("B18").Value = 0.016091653254452609
$ws.Range("A19").Value = -0.012091653267365832
$ws.Range("B19").Value = 0.01201692312626701
$ws.Range("A20").Value = -0.0080169231401221452
$ws.Range("B20").Value = 0.0080056490569671013
$ws.Range("A21").Value = -0.004005649070971451
$ws.Range("B21").Value = 0.0039999999858562063
$ws.Range("A22").Value = -0.118305864825901
$ws.Range("B22").Value = 0.11760459744187912
$ws.Range("A23").Value = -0.040498148239905163
$ws.Range("B23").Value = 0.040098903708696731
$ws.Range("A24").Value = -0.020098903778836608
$ws.Range("B24").Value = 0.019999999928890233
$ws.Range("A25").Value = -0.097215335362877298
$ws.Range("B25").Value = 0.097092006843778833
$ws.Range("A26").Value = -0.094592006866623279
$ws.Range("B26").Value = 0.094433589666856221
$ws.Range("A27").Value = -0.091933589690978756
$ws.Range("B27").Value = 0.090998140974291708
$ws.Range("A28").Value = -0.088998141002572417
$ws.Range("B28").Value = 0.088357613803444757
$ws.Range("A29").Value = -0.081357613850784993
$ws.Range("B29").Value = 0.081171701129505003
$ws.Range("A30").Value = -0.02117170133493218
$ws.Range("B30").Value = 0.021023995994000177
$ws.Range("A31").Value = -0.014023996044921105
$ws.Range("B31").Value = 0.01400141659699905
$ws.Range("A32").Value = -0.0040014166570099263
$ws.Range("B32").Value = 0.0039999999575659473

$ws.Columns.Item(1).ColumnWidth = 14.6
$ws.Columns.Item(2).ColumnWidth = 14.6
